$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.782872200012207
$ws.Range("B1").Value = 1.973571062088013
$ws.Range("C1").Value = 3.304929733276367
$ws.Range("D1").Value = 3.826452970504761
$ws.Range("E1").Value = 1.003474831581116
